$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("Raw")
$wsWh  = $wb.Worksheets.Item("Warehouse")

# ---------------------------------------------------------------------------
# Sheet "Raw": rows 18-25 lose their special (no-op) custom row formatting.
# These rows currently carry a "highlighted" style that renders the same as
# the default (fillId 0 even though applyFill=1), so strip it back to the
# plain / date look used elsewhere on the sheet.
# ---------------------------------------------------------------------------
18..25 | ForEach-Object { $wsRaw.Rows.Item($_).ClearFormats() }

# Re-apply the plain date style (same one used on B4:C4 etc.) to the date /
# WRR cells in those rows.
$wsRaw.Range("B4").Copy()
foreach ($r in 18,20,21,23,24,25) {
    $wsRaw.Range("B$r").PasteSpecial(-4122)
    $wsRaw.Range("C$r").PasteSpecial(-4122)
}
$wsRaw.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet "Warehouse"
# ---------------------------------------------------------------------------

# Row 31 (spatial-environment.R): the "waiting to hear from Josh on
# railroads" note is resolved -> clear the note and log a WRR run instead.
$wsWh.Rows.Item(31).ClearFormats()
$wsWh.Range("D31").ClearContents()
$wsWh.Range("B31").Value = 44925
$wsWh.Range("C31").Value = "WRR"
$wsWh.Range("B8").Copy()
$wsWh.Range("B31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 32 (spatial-other.R): log a WRR run, note stays the same.
$wsWh.Rows.Item(32).ClearFormats()
$wsWh.Range("B32").Value = 44925
$wsWh.Range("C32").Value = "WRR"
$wsWh.Range("B8").Copy()
$wsWh.Range("B32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 33 (spatial-parcel.R): log a WRR run, note stays the same.
$wsWh.Rows.Item(33).ClearFormats()
$wsWh.Range("B33").Value = 44931
$wsWh.Range("C33").Value = "WRR"
$wsWh.Range("B8").Copy()
$wsWh.Range("B33").PasteSpecial(-4122)
$wsWh.Range("C33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 35 (spatial-political.R): correct the last-successful-run date.
$wsWh.Range("B35").Value = 44566

# Row 36 (spatial-school.R): log a WRR run, highlighted like the other
# still-outstanding rows.
$wsWh.Range("B36").Value = 44936
$wsWh.Range("B36").NumberFormat = "mm-dd-yy"
$wsWh.Range("B36").Interior.ColorIndex = 6
$wsWh.Range("C36").Value = "WRR"
$wsWh.Range("A2").Copy()
$wsWh.Range("C36").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 37 (spatial-tax.R): log a (malformed) run date, highlighted like the
# other still-outstanding rows.
$wsWh.Range("B37").Value = "1/11/20203"
$wsWh.Range("B37").Interior.ColorIndex = 6
$wsWh.Range("B37").HorizontalAlignment = -4152
$wsWh.Range("C37").Value = "WRR"
$wsWh.Range("A2").Copy()
$wsWh.Range("C37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Warehouse sheet view: scroll down to the rows just touched and select A18.
# ---------------------------------------------------------------------------
$wsWh.Activate()
$wsWh.Range("A18").Select()
$excel.ActiveWindow.ScrollRow = 15
